$wb = $excel.ActiveWorkbook

# --- Sheet "summary" edits ---
$summary = $wb.Worksheets.Item("summary")

# Fill in the newly-scored query "5_wm" row (row 8) values that were blank
$summary.Range("F8").Value = 1
$summary.Range("J8").Value = 1
$summary.Range("N8").Value = 1
$summary.Range("R8").Value = 1
$summary.Range("V8").Value = 347

# Add the next two query rows (6_wm / 7_wm) against the same subset
$summary.Range("B9").Value = "6_wm"
$summary.Range("C9").Value = "cityofnewyork"
$summary.Range("B10").Value = "7_wm"
$summary.Range("C10").Value = "cityofnewyork"

# Update the active selection on the summary sheet
$summary.Range("V8").Select()

# --- Sheet "baseline" edits ---
$baseline = $wb.Worksheets.Item("baseline")
$baseline.Activate()
$baseline.Range("B5:C11").Select()

# Re-select the summary sheet as the active tab
$summary.Activate()
